$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (Total) sheet: insert a new leading row for 2022-Q4 and shift the
#    existing quarters down by one, recomputing the running index in col A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$totalRows = @(
    @(0, "2022-Q4", 30, 8.46),
    @(1, "2022-Q3", 26, 3.95),
    @(2, "2022-Q2", 6, 0.62),
    @(3, "2022-Q1", 3, 0.18),
    @(4, "2021-Q4", 1, 0.01),
    @(5, "2021-Q2", 7, 1.58),
    @(6, "2021-Q1", 19, 4.64),
    @(7, "2020-Q4", 7, 0.71)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Row 9 (2020-Q4) is brand new territory for this sheet; give its index cell
# (col A) the same look (bold / bordered / centered) as the rest of column A
# by copying the format from the cell right above it.
$total.Cells.Item(8, 1).Copy() | Out-Null
$total.Cells.Item(9, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q4" sheet itself, right after "总计" (i.e. right
#    before what is currently "2022-Q3"), with the quarterly fund holdings.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q4Data = @(
    @(0, '519714', '交银施罗德消费新驱动股票', '39.23', '88.65', '3.96', '1.5535', 10),
    @(1, '006604', '嘉实消费精选股票A', '25.39', '88.60', '4.86', '1.2340', 6),
    @(2, '005004', '交银施罗德品质升级混合A', '27.31', '88.42', '3.95', '1.0787', 10),
    @(3, '004868', '交银施罗德股息优化混合', '24.48', '90.05', '3.98', '0.9743', 10),
    @(4, '006605', '嘉实消费精选股票C', '15.66', '88.60', '4.86', '0.7611', 6),
    @(5, '014029', '浦银安盛红利精选混合C', '9.07', '71.06', '5.52', '0.5007', 1),
    @(6, '519115', '浦银安盛红利精选混合A', '7.60', '71.06', '5.52', '0.4195', 1),
    @(7, '013882', '交银施罗德品质升级混合C', '8.81', '88.42', '3.95', '0.3480', 10),
    @(8, '519710', '交银施罗德策略回报灵活配置混合', '7.64', '72.72', '3.76', '0.2873', 10),
    @(9, '001140', '工银总回报灵活配置混合A', '4.41', '81.48', '4.92', '0.2170', 3),
    @(10, '519125', '浦银安盛消费升级混合A', '2.11', '82.27', '8.23', '0.1737', 1),
    @(11, '013341', '工银核心机遇混合A', '5.83', '85.02', '2.88', '0.1679', 7),
    @(12, '013956', '华商医药消费精选混合A', '3.40', '80.90', '3.90', '0.1326', 9),
    @(13, '013293', '长城健康消费混合', '5.21', '86.80', '2.50', '0.1302', 10),
    @(14, '519176', '浦银安盛消费升级混合C', '0.92', '82.27', '8.23', '0.0757', 1),
    @(15, '014126', '华夏中证1000指数增强C', '8.66', '92.31', '0.80', '0.0693', 10),
    @(16, '159758', '华夏中证红利质量ETF', '1.81', '99.33', '3.67', '0.0664', 4),
    @(17, '002212', '嘉实新起航灵活配置混合A', '1.05', '79.36', '4.67', '0.0490', 9),
    @(18, '013342', '工银核心机遇混合C', '1.68', '85.02', '2.88', '0.0484', 7),
    @(19, '519991', '长信双利优选混合A', '1.11', '91.65', '4.11', '0.0456', 10),
    @(20, '006396', '长信双利优选混合E', '1.11', '91.65', '4.11', '0.0456', 10),
    @(21, '233015', '大摩量化配置混合A', '1.36', '93.08', '2.55', '0.0347', 9),
    @(22, '000649', '长城久鑫灵活配置混合A', '0.46', '90.33', '3.22', '0.0148', 3),
    @(23, '002512', '长城久润混合', '0.32', '93.69', '4.10', '0.0131', 4),
    @(24, '013957', '华商医药消费精选混合C', '0.26', '80.90', '3.90', '0.0101', 9),
    @(25, '014125', '华夏中证1000指数增强A', '1.04', '92.31', '0.80', '0.0083', 10),
    @(26, '011477', '工银总回报灵活配置混合C', '0.02', '81.48', '4.92', '0.0010', 3),
    @(27, '016264', '嘉实新起航灵活配置混合C', '0.01', '79.36', '4.67', '0.0005', 9),
    @(28, '008305', '大摩量化配置混合C', '0.01', '93.08', '2.55', '0.0003', 9),
    @(29, '017461', '长城久鑫灵活配置混合C', '0.00', '90.33', '3.22', '0', 3)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Last fund's market-value holding (G31) is genuinely zero, stored as a
# number rather than the "x.xxxx" text used everywhere else in that column.
$q4.Cells.Item(31, 7).Value = 0

# Header row + the running-index column share the workbook's one accent
# style: bold, thin box border, centered both ways.
$headerRange = $q4.Range($q4.Cells.Item(1, 2), $q4.Cells.Item(1, 8))
$indexRange = $q4.Range($q4.Cells.Item(2, 1), $q4.Cells.Item(31, 1))
foreach ($rng in @($headerRange, $indexRange)) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
    $rng.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
    $rng.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
}
